$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, shifting existing ticker rows down.
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "Ticker"

# Update the active selection to A2 (as in the saved workbook).
$ws.Range("A2").Select()

# The duplicate-value conditional formatting should cover only the data
# rows (A2:A44), not the new header row.
$fc = $ws.Range("A1:A43").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A2:A44"))
